$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Productos")

$ws.Range("C3").Value = 780
$ws.Range("C4").Value = 780
$ws.Range("C5").Value = 780
$ws.Range("C13").Value = 8
$ws.Range("C14").Value = 0.67
$ws.Range("C20").Value = 6.5

$ws.Range("C24").Select()
